# Recipient Transportation Mode Fractions.xlsx -- "test commit rmi-data branch"
#
# 1. Bump the "last updated" date on the About sheet.
# 2. Refresh a batch of battery-electric (col C) / diesel (col F) vehicle
#    counts on the raw "all_csv_SYVbT-passenger" data sheet with newer
#    source numbers, and blank out two cells (row 125 / Maine) that no
#    longer have source data.
# 3. Leave the UI in the state the author left it in: cursor parked on
#    F5 of "State RTMF", and "all_csv_SYVbT-passenger" as the selected /
#    fully-selected active tab (previously "RTMF-freight" was active).

$wb = $excel.ActiveWorkbook

# --- 1. About!C1 date bump (11/8/2022 -> 1/18/2023) ---------------------
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = 44944

# --- 2. Updated vehicle-count figures on all_csv_SYVbT-passenger --------
$ws = $wb.Worksheets.Item("all_csv_SYVbT-passenger")

$ws.Range("C17").Value = 3
$ws.Range("C23").Value = 24
$ws.Range("C29").Value = 450
$ws.Range("F29").Value = 27
$ws.Range("C35").Value = 72
$ws.Range("F35").Value = 0
$ws.Range("C41").Value = 0
$ws.Range("F41").Value = 5
$ws.Range("C47").Value = 0
$ws.Range("F47").Value = 0
$ws.Range("C53").Value = 29
$ws.Range("F53").Value = 11
$ws.Range("C59").Value = 41
$ws.Range("C83").Value = 279
$ws.Range("F83").Value = 21
$ws.Range("C89").Value = 17
$ws.Range("C107").Value = 31
$ws.Range("C113").Value = 166
$ws.Range("F113").Value = 43
$ws.Range("C119").Value = 50
$ws.Range("F119").Value = 5

# Row 125 (Maine): no longer has source data for these two modes -- clear
# the cells entirely rather than leaving a literal 0 behind.
$ws.Range("C125").ClearContents()
$ws.Range("F125").ClearContents()

$ws.Range("C131").Value = 5
$ws.Range("F131").Value = 0
$ws.Range("C137").Value = 26
$ws.Range("F137").Value = 1
$ws.Range("C143").Value = 32
$ws.Range("F143").Value = 0
$ws.Range("C161").Value = 20
$ws.Range("F161").Value = 0
$ws.Range("C185").Value = 190
$ws.Range("F185").Value = 25
$ws.Range("F191").Value = 7
$ws.Range("C203").Value = 835
$ws.Range("F203").Value = 17
$ws.Range("C209").Value = 24
$ws.Range("C215").Value = 4
$ws.Range("C221").Value = 72
$ws.Range("F221").Value = 1
$ws.Range("C227").Value = 297
$ws.Range("F227").Value = 0
$ws.Range("C233").Value = 0
$ws.Range("F233").Value = 0
$ws.Range("C251").Value = 2
$ws.Range("F251").Value = 4
$ws.Range("C257").Value = 94
$ws.Range("F257").Value = 19
$ws.Range("C263").Value = 37
$ws.Range("F263").Value = 6
$ws.Range("C269").Value = 19
$ws.Range("F269").Value = 19
$ws.Range("C281").Value = 40
$ws.Range("F281").Value = 2
$ws.Range("C287").Value = 4
$ws.Range("F287").Value = 0

# --- 3. Restore the author's on-screen cursor / active-sheet state ------
# "State RTMF" keeps a remembered cursor on F5 even though it isn't the
# sheet left on top when the file was saved.
$wsState = $wb.Worksheets.Item("State RTMF")
$wsState.Activate()
$wsState.Range("F5").Select()

# The sheet actually left active/selected (whole sheet selected) is
# "all_csv_SYVbT-passenger" -- do this last so it ends up the active tab.
$ws.Activate()
$ws.Cells.Select()
